$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'silicone knee pads',
    'skins compression men',
    'sliding knee pads baseball',
    'small football knee pads',
    'snowboarding knee pad',
    'soccer knee pad',
    'soccer knee protector',
    'soccer pants for boys',
    'spandex pants for men',
    'spats bjj',
    'sports basketball pants',
    'sports leggings men',
    'squat protector',
    'strip pants men',
    'tactical pants with knee pads',
    'thermal baselayer men',
    'thermal compression pants men',
    'thermal leggings men',
    'thermal winter pants',
    'thin knee pads',
    'tight pant',
    'tights black',
    'tights for men',
    'under amour fleece leggings',
    'under armour leggings',
    'underarmor thermal pants mens',
    'volleyball clothes men',
    'volleyball compression knee pads',
    'volleyball knee pads adult',
    'volleyball knee pads youth',
    'warm compression pants mens',
    'weightlifting equipment',
    'white basketball knee pads',
    'white tights mens basketball',
    'womens basketball knee pads',
    'workout pads for hands',
    'workout squat pad',
    'wrestling knee pads adult',
    'wrestling tights boys',
    'xl volleyball knee pads',
    'youth basketball knee',
    'youth basketball leggings for boys',
    'youth basketball pants girls',
    'youth football knee pads',
    'youth knee pads for skating',
    'youth knee pads mountain bike',
    'youth soccer gear',
    'youth sports tights',
    'mens running tights capri',
    'mens basketball pants tall',
    'volleyball hip protectors',
    'knee compression basketball',
    'sports compression pants',
    'basketball leggings youth',
    'knee pads for running',
    'wrestling pants for men',
    'male compression pants',
    'cheap compression pants men',
    'men compression pants pack',
    'knee pad men',
    'volleyball knee pads for youth',
    'leggings for man',
    'volleyball knee pads xl',
    'spandex leggings men',
    'knees compression',
    'athletic boys pants',
    'soccer knee pads',
    'adults knee pads',
    'leg compression leggings',
    'big and tall tights for men',
    'wrestling mens apparel',
    'mens hiking leggings',
    'softball pants mens',
    'running compression knee',
    'compression pads for basketball',
    'gym tights',
    'boys soccer leggings',
    'hockey pants men',
    'knee pads for volleyball for men',
    'baseball pants men',
    'mens running pants',
    'sports pants for men',
    'boys snowboarding pants',
    'clothing protector',
    'baseball pants for youth',
    'mens gym pants',
    'football pants adult',
    'knee pads sports',
    'protective knee pads',
    'girls volleyball pads',
    'pants youth',
    'big boys compression leggings',
    'basketball spandex',
    'mens basketball clothing',
    'boys tights for sports',
    'leggings for cycling',
    'hiking capri men',
    'baseball pants men black',
    'football tights youth',
    'mens 3/4 tights'
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

